$d = $word.ActiveDocument

$d.Content.Find.Execute('389×7=2723', $true, $false, $false, $false, $false, $true, 1, $false, '882×2=1764', 2) | Out-Null
$d.Content.Find.Execute('266×8=2128', $true, $false, $false, $false, $false, $true, 1, $false, '149×5=745', 2) | Out-Null
$d.Content.Find.Execute('384×9=3456', $true, $false, $false, $false, $false, $true, 1, $false, '437×5=2185', 2) | Out-Null
$d.Content.Find.Execute('163×8=1304', $true, $false, $false, $false, $false, $true, 1, $false, '641×4=2564', 2) | Out-Null
$d.Content.Find.Execute('734×5=3670', $true, $false, $false, $false, $false, $true, 1, $false, '919×4=3676', 2) | Out-Null
$d.Content.Find.Execute('154×6=924', $true, $false, $false, $false, $false, $true, 1, $false, '653×4=2612', 2) | Out-Null
$d.Content.Find.Execute('323×9=2907', $true, $false, $false, $false, $false, $true, 1, $false, '916×2=1832', 2) | Out-Null
$d.Content.Find.Execute('234×4=936', $true, $false, $false, $false, $false, $true, 1, $false, '495×3=1485', 2) | Out-Null
$d.Content.Find.Execute('810×5=4050', $true, $false, $false, $false, $false, $true, 1, $false, '321×3=963', 2) | Out-Null
$d.Content.Find.Execute('401×9=3609', $true, $false, $false, $false, $false, $true, 1, $false, '213×7=1491', 2) | Out-Null
$d.Content.Find.Execute('977×9=8793', $true, $false, $false, $false, $false, $true, 1, $false, '666×3=1998', 2) | Out-Null
$d.Content.Find.Execute('365×5=1825', $true, $false, $false, $false, $false, $true, 1, $false, '628×8=5024', 2) | Out-Null
$d.Content.Find.Execute('693×6=4158', $true, $false, $false, $false, $false, $true, 1, $false, '652×7=4564', 2) | Out-Null
$d.Content.Find.Execute('636×2=1272', $true, $false, $false, $false, $false, $true, 1, $false, '958×5=4790', 2) | Out-Null
$d.Content.Find.Execute('280×3=840', $true, $false, $false, $false, $false, $true, 1, $false, '442×8=3536', 2) | Out-Null
$d.Content.Find.Execute('458×6=2748', $true, $false, $false, $false, $false, $true, 1, $false, '376×4=1504', 2) | Out-Null
$d.Content.Find.Execute('133×7=931', $true, $false, $false, $false, $false, $true, 1, $false, '973×9=8757', 2) | Out-Null
$d.Content.Find.Execute('867×5=4335', $true, $false, $false, $false, $false, $true, 1, $false, '360×3=1080', 2) | Out-Null
$d.Content.Find.Execute('984×9=8856', $true, $false, $false, $false, $false, $true, 1, $false, '490×5=2450', 2) | Out-Null
$d.Content.Find.Execute('205×8=1640', $true, $false, $false, $false, $false, $true, 1, $false, '233×5=1165', 2) | Out-Null
$d.Content.Find.Execute('495×5=2475', $true, $false, $false, $false, $false, $true, 1, $false, '939×5=4695', 2) | Out-Null
$d.Content.Find.Execute('939×7=6573', $true, $false, $false, $false, $false, $true, 1, $false, '551×9=4959', 2) | Out-Null
$d.Content.Find.Execute('757×2=1514', $true, $false, $false, $false, $false, $true, 1, $false, '932×2=1864', 2) | Out-Null
$d.Content.Find.Execute('394×5=1970', $true, $false, $false, $false, $false, $true, 1, $false, '386×9=3474', 2) | Out-Null
$d.Content.Find.Execute('867×6=5202', $true, $false, $false, $false, $false, $true, 1, $false, '111×3=333', 2) | Out-Null
